$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new label value for column F ("labels").
# These normalize label casing (and, for multi-valued "||"-joined
# labels, reorder the individual tags) per the corrected annotations.
$updates = @{
    4   = "135_product_information"
    7   = "93_referral_statement"
    15  = "18_hazards_to_humans_and_domestic_animals"
    16  = "ppe"
    17  = "ppe"
    18  = "ppe"
    22  = "application instructions"
    23  = "env warning - species"
    24  = "env warning - species"
    25  = "env warning - water"
    27  = "env warning - water || off target movement"
    28  = "env warning - water"
    48  = "application instructions"
    49  = "application instructions"
    51  = "application instructions"
    52  = "application instructions"
    54  = "134_non-agriculture_use_requirements"
    57  = "135_product_information"
    58  = "135_product_information"
    60  = "use restrictions"
    61  = "use restrictions"
    62  = "use restrictions"
    63  = "off target movement"
    70  = "use restrictions || off target movement"
    76  = "safety procedures || application instructions"
    77  = "mixing"
    78  = "mixing"
    79  = "mixing"
    80  = "mixing"
    81  = "mixing"
    84  = "application instructions"
    85  = "application instructions"
    86  = "application instructions"
    87  = "use restrictions || application instructions"
    88  = "use restrictions || application instructions"
    89  = "irrigation || application instructions || chemigation"
    90  = "safety procedures || use restrictions"
    91  = "use restrictions || irrigation"
    92  = "use restrictions || irrigation"
    93  = "chemigation"
    95  = "application instructions"
    96  = "application instructions"
    118 = "use restrictions"
    198 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
